{"js": "// Replace the two-digit multiplication problems in the worksheet table\n// with the newly generated set of problems. Each old expression is\n// unique within the document, so a literal search-and-replace on each\n// pair is sufficient and order-independent.\nconst replacements = [\n  [\"41\u00d768=\", \"57\u00d796=\"],\n  [\"24\u00d746=\", \"11\u00d739=\"],\n  [\"32\u00d756=\", \"59\u00d774=\"],\n  [\"63\u00d765=\", \"31\u00d790=\"],\n  [\"41\u00d774=\", \"34\u00d794=\"],\n  [\"31\u00d715=\", \"35\u00d794=\"],\n  [\"45\u00d770=\", \"81\u00d754=\"],\n  [\"87\u00d783=\", \"19\u00d747=\"],\n  [\"50\u00d722=\", \"73\u00d718=\"],\n  [\"76\u00d789=\", \"68\u00d742=\"],\n  [\"47\u00d730=\", \"32\u00d750=\"],\n  [\"22\u00d732=\", \"12\u00d731=\"],\n  [\"88\u00d721=\", \"91\u00d774=\"],\n  [\"64\u00d726=\", \"57\u00d799=\"],\n  [\"19\u00d721=\", \"45\u00d746=\"],\n  [\"59\u00d756=\", \"87\u00d762=\"],\n  [\"97\u00d732=\", \"39\u00d774=\"],\n  [\"51\u00d793=\", \"43\u00d789=\"],\n  [\"18\u00d716=\", \"83\u00d796=\"],\n  [\"62\u00d791=\", \"68\u00d756=\"],\n  [\"67\u00d767=\", \"14\u00d741=\"],\n  [\"87\u00d756=\", \"58\u00d719=\"],\n  [\"67\u00d781=\", \"41\u00d730=\"],\n  [\"15\u00d749=\", \"43\u00d715=\"],\n  [\"73\u00d724=\", \"16\u00d751=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit multiplication problems in the worksheet table\n# with the newly generated set of problems. Each old expression is\n# unique within the document, so a literal Find/Replace on each pair\n# (scoped to the whole document) is sufficient and order-independent.\n\n$replacements = @(\n    @(\"41\u00d768=\", \"57\u00d796=\"),\n    @(\"24\u00d746=\", \"11\u00d739=\"),\n    @(\"32\u00d756=\", \"59\u00d774=\"),\n    @(\"63\u00d765=\", \"31\u00d790=\"),\n    @(\"41\u00d774=\", \"34\u00d794=\"),\n    @(\"31\u00d715=\", \"35\u00d794=\"),\n    @(\"45\u00d770=\", \"81\u00d754=\"),\n    @(\"87\u00d783=\", \"19\u00d747=\"),\n    @(\"50\u00d722=\", \"73\u00d718=\"),\n    @(\"76\u00d789=\", \"68\u00d742=\"),\n    @(\"47\u00d730=\", \"32\u00d750=\"),\n    @(\"22\u00d732=\", \"12\u00d731=\"),\n    @(\"88\u00d721=\", \"91\u00d774=\"),\n    @(\"64\u00d726=\", \"57\u00d799=\"),\n    @(\"19\u00d721=\", \"45\u00d746=\"),\n    @(\"59\u00d756=\", \"87\u00d762=\"),\n    @(\"97\u00d732=\", \"39\u00d774=\"),\n    @(\"51\u00d793=\", \"43\u00d789=\"),\n    @(\"18\u00d716=\", \"83\u00d796=\"),\n    @(\"62\u00d791=\", \"68\u00d756=\"),\n    @(\"67\u00d767=\", \"14\u00d741=\"),\n    @(\"87\u00d756=\", \"58\u00d719=\"),\n    @(\"67\u00d781=\", \"41\u00d730=\"),\n    @(\"15\u00d749=\", \"43\u00d715=\"),\n    @(\"73\u00d724=\", \"16\u00d751=\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n\n    $find.Execute(\n        $oldText,   # FindText\n        $true,      # MatchCase\n        $false,     # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap (wdFindContinue)\n        $false,     # Format\n        $newText,   # ReplaceWith\n        2           # Replace (wdReplaceAll)\n    )\n}\n"}
